$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (GitHub Actions scheduled update).
# Price cells (column D) are stored as text in this sheet (e.g. "27.100.13",
# "1.006") rather than numbers, so values that look numeric are written with
# a leading apostrophe to force text entry and avoid Excel silently
# reinterpreting/reformatting them as numbers (which would drop trailing
# zeros, e.g. "0.09000" -> 0.09, or merge "27.123.00" into a float).

$ws.Range("D2").Value = "27.123.00"
$ws.Range("E2").Value = "  -1.93%  "

$ws.Range("D3").Value = "1.804.73"
$ws.Range("E3").Value = "  -2.18%  "

$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").Value = "'309.21"
$ws.Range("E5").Value = "  -1.91%  "

$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  +0.20%  "

$ws.Range("D8").Value = "'0.3617"
$ws.Range("E8").Value = "  -2.39%  "

$ws.Range("E9").Value = "  -1.40%  "

$ws.Range("D10").Value = "'0.8460"
$ws.Range("E10").Value = "  -3.74%  "

$ws.Range("D11").Value = "'20.35"
$ws.Range("E11").Value = "  -3.04%  "

$ws.Range("D12").Value = "1.783.52"
$ws.Range("E12").Value = "  -1.64%  "

$ws.Range("D13").Value = "'5.305"
$ws.Range("E13").Value = "  -3.25%  "

$ws.Range("D14").Value = "'6.394"
$ws.Range("E14").Value = "  -3.23%  "

$ws.Range("D15").Value = "'0.06798"
$ws.Range("E15").Value = "  -2.27%  "

$ws.Range("D16").Value = "'1.008"
$ws.Range("E16").Value = "  +0.42%  "

$ws.Range("D17").Value = "'81.24"
$ws.Range("E17").Value = "  +0.01%  "

$ws.Range("D18").Value = "'0.000008771"
$ws.Range("E18").Value = "  -2.96%  "

$ws.Range("E19").Value = "  +0.40%  "

$ws.Range("D20").Value = "'15.03"
$ws.Range("E20").Value = "  -3.67%  "

$ws.Range("D21").Value = "27.218.09"
$ws.Range("E21").Value = "  -1.24%  "

$ws.Range("D22").Value = "'5.093"
$ws.Range("E22").Value = "  -0.91%  "

$ws.Range("D23").Value = "'11.10"
$ws.Range("E23").Value = "  +0.82%  "

$ws.Range("D24").Value = "2.053.03"
$ws.Range("E24").Value = "  -4.18%  "

$ws.Range("D25").Value = "'1.953"
$ws.Range("E25").Value = "  -1.93%  "

$ws.Range("D26").Value = "'153.30"
$ws.Range("E26").Value = "  -0.41%  "

$ws.Range("E27").Value = "  -3.59%  "

$ws.Range("D28").Value = "'5.040"
$ws.Range("E28").Value = "  -5.18%  "

$ws.Range("D29").Value = "'113.96"
$ws.Range("E29").Value = "  -1.92%  "

$ws.Range("D30").Value = "'1.663"
$ws.Range("E30").Value = "  -11.33%  "

$ws.Range("D31").Value = "'0.09000"
$ws.Range("E31").Value = "  +0.68%  "

$ws.Range("D32").Value = "'0.7367"
$ws.Range("E32").Value = "  -6.44%  "

$ws.Range("D33").Value = "'2.883"
$ws.Range("E33").Value = "  -3.09%  "

$ws.Range("D34").Value = "'4.378"
$ws.Range("E34").Value = "  -5.22%  "

$ws.Range("E35").Value = "  -6.61%  "

$ws.Range("E36").Value = "  +0.25%  "

$ws.Range("D37").Value = "'1.081"
$ws.Range("E37").Value = "  -2.02%  "

$ws.Range("D38").Value = "'0.05161"
$ws.Range("E38").Value = "  -5.21%  "

$ws.Range("D39").Value = "'0.01912"
$ws.Range("E39").Value = "  -2.70%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.4999"
$ws.Range("E40").Value = "  -3.50%  "

$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "'0.1636"
$ws.Range("E41").Value = "  -3.25%  "

$ws.Range("D42").Value = "'2.615"
$ws.Range("E42").Value = "  -8.33%  "

$ws.Range("D43").Value = "'8.138"
$ws.Range("E43").Value = "  -5.83%  "

$ws.Range("D44").Value = "'5.964"
$ws.Range("E44").Value = "  -12.21%  "

$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'105.17"
$ws.Range("E45").Value = "  -1.45%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'10.23"
$ws.Range("E46").Value = "  -3.84%  "

$ws.Range("E47").Value = "  +0.20%  "

$ws.Range("D48").Value = "'0.06334"
$ws.Range("E48").Value = "  -3.40%  "

$ws.Range("D49").Value = "'0.4550"
$ws.Range("E49").Value = "  -5.06%  "

$ws.Range("D50").Value = "'1.609"
$ws.Range("E50").Value = "  -3.44%  "

$ws.Range("D51").Value = "'1.723"
$ws.Range("E51").Value = "  -6.44%  "
